$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 5
$ws.Range("F2").Value = 217
$ws.Range("H2").Value = "bedrooms"
$ws.Range("I2").Value = "distractor"
$ws.Range("K2").Value = "f"
$ws.Range("L2").Value = "stimuli/img_scrdm.png"
$ws.Range("M2").Value = 78.675
$ws.Range("N2").Value = 57.9
$ws.Range("O2").Value = 68.2875
$ws.Range("P2").Value = 40
$ws.Range("Q2").Value = 7
$ws.Range("R2").Value = 7
$ws.Range("S2").Value = 7

# Row 3
$ws.Range("C3").Value = 5
$ws.Range("F3").Value = 218
$ws.Range("H3").Value = "kitchens"
$ws.Range("I3").Value = "target"
$ws.Range("K3").Value = "j"
$ws.Range("L3").Value = "stimuli/img_ce9vx.png"
$ws.Range("M3").Value = 75.9090909090909
$ws.Range("N3").Value = 57.12121212121212
$ws.Range("O3").Value = 66.51515151515152
$ws.Range("P3").Value = 33
$ws.Range("Q3").Value = 7
$ws.Range("R3").Value = 7
$ws.Range("S3").Value = 7

# Row 4
$ws.Range("C4").Value = 5
$ws.Range("F4").Value = 219
$ws.Range("H4").Value = "kitchens"
$ws.Range("I4").Value = "target"
$ws.Range("K4").Value = "j"
$ws.Range("L4").Value = "stimuli/img_p3hpc.png"
$ws.Range("M4").Value = 72.83333333333333
$ws.Range("N4").Value = 52.22222222222222
$ws.Range("O4").Value = 62.52777777777777
$ws.Range("P4").Value = 36
$ws.Range("Q4").Value = 6
$ws.Range("R4").Value = 6
$ws.Range("S4").Value = 6

# Row 5
$ws.Range("C5").Value = 5
$ws.Range("F5").Value = 220
$ws.Range("L5").Value = "stimuli/img_eatdk.png"
$ws.Range("N5").Value = 61.375
$ws.Range("O5").Value = 71.390625

# Row 6
$ws.Range("C6").Value = 5
$ws.Range("F6").Value = 221
$ws.Range("L6").Value = "stimuli/img_aplao.png"
$ws.Range("M6").Value = 64.0909090909091
$ws.Range("N6").Value = 40.75757575757576
$ws.Range("O6").Value = 52.42424242424242
$ws.Range("P6").Value = 33
$ws.Range("Q6").Value = 3
$ws.Range("R6").Value = 3
$ws.Range("S6").Value = 3

# Row 7
$ws.Range("C7").Value = 5
$ws.Range("F7").Value = 222
$ws.Range("L7").Value = "stimuli/img_iyxnj.png"
$ws.Range("M7").Value = 75.30555555555556
$ws.Range("N7").Value = 54.33333333333334
$ws.Range("O7").Value = 64.81944444444444
$ws.Range("P7").Value = 36
$ws.Range("Q7").Value = 6
$ws.Range("R7").Value = 6
$ws.Range("S7").Value = 6

# Row 8
$ws.Range("C8").Value = 5
$ws.Range("F8").Value = 223
$ws.Range("H8").Value = "kitchens"
$ws.Range("I8").Value = "target"
$ws.Range("K8").Value = "j"
$ws.Range("L8").Value = "stimuli/img_nyv2b.png"
$ws.Range("M8").Value = 11.91176470588235
$ws.Range("N8").Value = 6.852941176470588
$ws.Range("O8").Value = 9.382352941176471
$ws.Range("P8").Value = 34
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = 1
$ws.Range("S8").Value = 1

# Row 9
$ws.Range("C9").Value = 5
$ws.Range("F9").Value = 224
$ws.Range("L9").Value = "stimuli/img_cnyac.png"
$ws.Range("M9").Value = 69.1470588235294
$ws.Range("N9").Value = 47.8235294117647
$ws.Range("O9").Value = 58.48529411764706
$ws.Range("P9").Value = 34
$ws.Range("Q9").Value = 5
$ws.Range("R9").Value = 5
$ws.Range("S9").Value = 5

# Row 10
$ws.Range("C10").Value = 5
$ws.Range("F10").Value = 225
$ws.Range("H10").Value = "bedrooms"
$ws.Range("I10").Value = "distractor"
$ws.Range("K10").Value = "f"
$ws.Range("L10").Value = "stimuli/img_kn0we.png"
$ws.Range("M10").Value = 80.1590909090909
$ws.Range("N10").Value = 56.68181818181818
$ws.Range("O10").Value = 68.42045454545455
$ws.Range("P10").Value = 44
$ws.Range("Q10").Value = 7
$ws.Range("R10").Value = 7
$ws.Range("S10").Value = 7

# Row 11
$ws.Range("C11").Value = 5
$ws.Range("F11").Value = 226
$ws.Range("H11").Value = "living_rooms"
$ws.Range("I11").Value = "distractor"
$ws.Range("K11").Value = "f"
$ws.Range("L11").Value = "stimuli/img_g13d5.png"
$ws.Range("M11").Value = 73
$ws.Range("N11").Value = 51.51111111111111
$ws.Range("O11").Value = 62.25555555555556
$ws.Range("P11").Value = 45
$ws.Range("Q11").Value = 6
$ws.Range("R11").Value = 6
$ws.Range("S11").Value = 6

# Row 12
$ws.Range("C12").Value = 5
$ws.Range("F12").Value = 227
$ws.Range("L12").Value = "stimuli/img_3gm8h.png"
$ws.Range("M12").Value = 65.07894736842105
$ws.Range("N12").Value = 43.92105263157895
$ws.Range("O12").Value = 54.5
$ws.Range("P12").Value = 38
$ws.Range("Q12").Value = 4
$ws.Range("R12").Value = 4
$ws.Range("S12").Value = 4

# Row 13
$ws.Range("C13").Value = 5
$ws.Range("F13").Value = 228
$ws.Range("L13").Value = "stimuli/img_yeh72.png"
$ws.Range("M13").Value = 68.66666666666667
$ws.Range("N13").Value = 45.21212121212121
$ws.Range("O13").Value = 56.93939393939394
$ws.Range("P13").Value = 33
$ws.Range("Q13").Value = 4
$ws.Range("R13").Value = 4
$ws.Range("S13").Value = 4

# Row 14
$ws.Range("C14").Value = 5
$ws.Range("F14").Value = 229
$ws.Range("L14").Value = "stimuli/img_uwv6y.png"
$ws.Range("M14").Value = 78.88888888888889
$ws.Range("N14").Value = 59.30555555555556
$ws.Range("O14").Value = 69.09722222222223
$ws.Range("P14").Value = 36
$ws.Range("Q14").Value = 8
$ws.Range("R14").Value = 8
$ws.Range("S14").Value = 8

# Row 15
$ws.Range("C15").Value = 5
$ws.Range("F15").Value = 230
$ws.Range("L15").Value = "stimuli/img_a8wvq.png"
$ws.Range("M15").Value = 86.25925925925925
$ws.Range("N15").Value = 66.25925925925925
$ws.Range("O15").Value = 76.25925925925925
$ws.Range("P15").Value = 27
$ws.Range("Q15").Value = 10
$ws.Range("R15").Value = 10
$ws.Range("S15").Value = 10

# Row 16
$ws.Range("C16").Value = 5
$ws.Range("F16").Value = 231
$ws.Range("H16").Value = "living_rooms"
$ws.Range("I16").Value = "distractor"
$ws.Range("K16").Value = "f"
$ws.Range("L16").Value = "stimuli/img_x9w7o.png"
$ws.Range("M16").Value = 92.38888888888889
$ws.Range("N16").Value = 72.94444444444444
$ws.Range("O16").Value = 82.66666666666666
$ws.Range("P16").Value = 36
$ws.Range("Q16").Value = 10
$ws.Range("R16").Value = 10
$ws.Range("S16").Value = 10

# Row 17
$ws.Range("C17").Value = 5
$ws.Range("F17").Value = 232
$ws.Range("L17").Value = "stimuli/img_463mq.png"
$ws.Range("M17").Value = 51.35294117647059
$ws.Range("N17").Value = 30.20588235294118
$ws.Range("O17").Value = 40.77941176470588
$ws.Range("Q17").Value = 2
$ws.Range("R17").Value = 2
$ws.Range("S17").Value = 2

# Row 18
$ws.Range("C18").Value = 5
$ws.Range("F18").Value = 233
$ws.Range("H18").Value = "kitchens"
$ws.Range("I18").Value = "target"
$ws.Range("K18").Value = "j"
$ws.Range("L18").Value = "stimuli/img_9mky8.png"
$ws.Range("M18").Value = 84.32352941176471
$ws.Range("N18").Value = 65.17647058823529
$ws.Range("O18").Value = 74.75
$ws.Range("P18").Value = 34
$ws.Range("Q18").Value = 9
$ws.Range("R18").Value = 9
$ws.Range("S18").Value = 9

# Row 19
$ws.Range("C19").Value = 5
$ws.Range("F19").Value = 234
$ws.Range("L19").Value = "stimuli/img_inqod.png"
$ws.Range("M19").Value = 70.84848484848484
$ws.Range("N19").Value = 50.63636363636363
$ws.Range("O19").Value = 60.74242424242424
$ws.Range("P19").Value = 33
$ws.Range("Q19").Value = 5
$ws.Range("R19").Value = 5
$ws.Range("S19").Value = 5

# Row 20
$ws.Range("C20").Value = 5
$ws.Range("F20").Value = 235
$ws.Range("L20").Value = "stimuli/img_ye5sl.png"
$ws.Range("M20").Value = 53.2258064516129
$ws.Range("N20").Value = 34.45161290322581
$ws.Range("O20").Value = 43.83870967741936
$ws.Range("P20").Value = 31
$ws.Range("Q20").Value = 2
$ws.Range("R20").Value = 2
$ws.Range("S20").Value = 2

# Row 21
$ws.Range("C21").Value = 5
$ws.Range("F21").Value = 236
$ws.Range("L21").Value = "stimuli/img_wyl6z.png"
$ws.Range("M21").Value = 59.8235294117647
$ws.Range("N21").Value = 36.23529411764706
$ws.Range("O21").Value = 48.02941176470588
$ws.Range("P21").Value = 34
$ws.Range("Q21").Value = 3
$ws.Range("R21").Value = 3
$ws.Range("S21").Value = 3

# Row 22
$ws.Range("C22").Value = 5
$ws.Range("F22").Value = 237
$ws.Range("L22").Value = "stimuli/img_t90e2.png"
$ws.Range("M22").Value = 83.0625
$ws.Range("N22").Value = 61.96875
$ws.Range("O22").Value = 72.515625
$ws.Range("Q22").Value = 9
$ws.Range("R22").Value = 9
$ws.Range("S22").Value = 9

# Row 23
$ws.Range("C23").Value = 5
$ws.Range("F23").Value = 238
$ws.Range("L23").Value = "stimuli/img_60242.png"
$ws.Range("M23").Value = 78.33333333333333
$ws.Range("N23").Value = 57.57575757575758
$ws.Range("O23").Value = 67.95454545454545
$ws.Range("Q23").Value = 7
$ws.Range("R23").Value = 7
$ws.Range("S23").Value = 7

# Row 24
$ws.Range("C24").Value = 5
$ws.Range("F24").Value = 239
$ws.Range("L24").Value = "stimuli/img_jpjeg.png"
$ws.Range("M24").Value = 90.90697674418605
$ws.Range("N24").Value = 74.3953488372093
$ws.Range("O24").Value = 82.65116279069767
$ws.Range("P24").Value = 43
$ws.Range("Q24").Value = 10
$ws.Range("R24").Value = 10
$ws.Range("S24").Value = 10

# Row 25
$ws.Range("C25").Value = 5
$ws.Range("F25").Value = 240
$ws.Range("H25").Value = "living_rooms"
$ws.Range("I25").Value = "distractor"
$ws.Range("K25").Value = "f"
$ws.Range("L25").Value = "stimuli/img_53nbn.png"
$ws.Range("M25").Value = 73.28888888888889
$ws.Range("N25").Value = 51.15555555555556
$ws.Range("O25").Value = 62.22222222222223
$ws.Range("P25").Value = 45
$ws.Range("Q25").Value = 6
$ws.Range("R25").Value = 6
$ws.Range("S25").Value = 6

# Row 26
$ws.Range("C26").Value = 5
$ws.Range("F26").Value = 241
$ws.Range("L26").Value = "stimuli/img_7wul8.png"
$ws.Range("M26").Value = 43.03030303030303
$ws.Range("N26").Value = 25.54545454545455
$ws.Range("O26").Value = 34.28787878787879
$ws.Range("Q26").Value = 1
$ws.Range("R26").Value = 1
$ws.Range("S26").Value = 1

# Row 27
$ws.Range("C27").Value = 5
$ws.Range("F27").Value = 242
$ws.Range("H27").Value = "kitchens"
$ws.Range("I27").Value = "target"
$ws.Range("K27").Value = "j"
$ws.Range("L27").Value = "stimuli/img_d8xbu.png"
$ws.Range("M27").Value = 91.36363636363636
$ws.Range("N27").Value = 73.18181818181819
$ws.Range("O27").Value = 82.27272727272728
$ws.Range("P27").Value = 33
